# "changes in Appium Capabilities"
# Flip the two test-result statuses in the Status column (D) of the
# UserDataForm sheet: row 2 goes from FAIL to PASS, row 5 goes from
# PASS to FAIL. Also update the sheet's saved selection to D2:D5 to
# match where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "PASS"
$ws.Range("D5").Value = "FAIL"

$ws.Range("D2:D5").Select() | Out-Null
